# Update the Cascades sheet so each stage references the corresponding
# "all_*" Characteristic code name (for which data is actually provided
# in the databook) instead of an ad-hoc comma-separated list of
# compartments.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cascades")

$ws.Range("B2").Value = "all_people"
$ws.Range("B3").Value = "all_dx"
$ws.Range("B4").Value = "all_ever_linked"
$ws.Range("B5").Value = "all_curr_linked"
$ws.Range("B6").Value = "all_tx"
$ws.Range("B7").Value = "all_vs"

# Reflect the reviewer's selection while cross-checking the Characteristics
# codes used above, then return focus/selection to the Cascades sheet
# (which stays the active tab).
$wsChar = $wb.Worksheets.Item("Characteristics")
$wsChar.Range("A2:A7").Select()

$ws.Range("B2:B7").Select()
